# Generate Report for Handoff
#
# The localization CI re-ran and produced a fresh handoff package: the
# per-locale status flips from "In Translation" to "Ready for handoff" and
# the associated timestamps advance a few seconds. The Status/"Latest HO
# Xliff Generate Date"/"Latest Handoff Datetime" columns are widened on all
# three sheets so the longer "Ready for handoff" label isn't clipped.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ColumnWidth (characters) this host rounds to whole-pixel granularity on
# save; 16.33 is the input that lands closest to the target stored width.
$newColWidth = 16.33

# --- Overview sheet ------------------------------------------------------
# E2 / F2: per-locale status ("In Translation" -> "Ready for handoff")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
# G2: Latest HO Xliff Generate Date
$overview.Range("G2").Value = "2016-09-05 07:08:35"

# Widen columns E and F (zh-cn / de-de status) to fit the longer text
$overview.Columns.Item(5).ColumnWidth = $newColWidth
$overview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet -----------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-05 07:08:31"
$zhcn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet -----------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-05 07:08:35"
$dede.Columns.Item(3).ColumnWidth = $newColWidth
